{"js": "const body = context.document.body;\n\n// Map of old exercise text -> new exercise text.\nconst replacements = [\n  [\"19\u00f74=4, 3\", \"54\u00f77=7, 5\"],\n  [\"23\u00f75=4, 3\", \"74\u00f77=10, 4\"],\n  [\"96\u00f78=12, 0\", \"22\u00f72=11, 0\"],\n  [\"83\u00f76=13, 5\", \"79\u00f73=26, 1\"],\n  [\"35\u00f77=5, 0\", \"83\u00f77=11, 6\"],\n  [\"66\u00f75=13, 1\", \"65\u00f72=32, 1\"],\n  [\"44\u00f79=4, 8\", \"91\u00f75=18, 1\"],\n  [\"94\u00f76=15, 4\", \"96\u00f77=13, 5\"],\n  [\"82\u00f76=13, 4\", \"85\u00f73=28, 1\"],\n  [\"87\u00f79=9, 6\", \"67\u00f79=7, 4\"],\n  [\"11\u00f74=2, 3\", \"25\u00f72=12, 1\"],\n  [\"38\u00f79=4, 2\", \"56\u00f72=28, 0\"],\n  [\"33\u00f74=8, 1\", \"25\u00f74=6, 1\"],\n  [\"23\u00f79=2, 5\", \"71\u00f74=17, 3\"],\n  [\"12\u00f79=1, 3\", \"58\u00f77=8, 2\"],\n  [\"44\u00f77=6, 2\", \"54\u00f73=18, 0\"],\n  [\"81\u00f77=11, 4\", \"41\u00f73=13, 2\"],\n  [\"74\u00f77=10, 4\", \"80\u00f75=16, 0\"],\n  [\"55\u00f78=6, 7\", \"62\u00f76=10, 2\"],\n  [\"32\u00f79=3, 5\", \"97\u00f75=19, 2\"],\n  [\"37\u00f77=5, 2\", \"69\u00f79=7, 6\"],\n  [\"32\u00f75=6, 2\", \"72\u00f74=18, 0\"],\n  [\"50\u00f73=16, 2\", \"76\u00f76=12, 4\"],\n  [\"38\u00f75=7, 3\", \"89\u00f79=9, 8\"],\n  [\"84\u00f73=28, 0\", \"14\u00f78=1, 6\"],\n];\n\n// Phase 1: swap each old value for a unique placeholder so that\n// a later replacement can never accidentally match text that was\n// just inserted by an earlier one (some new values equal other old values).\nconst placeholders = replacements.map((_, i) => `@@PLACEHOLDER${i}@@`);\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText] = replacements[i];\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load('items');\n  await context.sync();\n  for (const r of found.items) {\n    r.insertText(placeholders[i], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Phase 2: swap each placeholder for its final value.\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const found = body.search(placeholders[i], { matchCase: true, matchWholeWord: false });\n  found.load('items');\n  await context.sync();\n  for (const r of found.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Old exercise text -> new exercise text (in document order).\n$replacements = @(\n    ,@(\"19\u00f74=4, 3\", \"54\u00f77=7, 5\")\n    ,@(\"23\u00f75=4, 3\", \"74\u00f77=10, 4\")\n    ,@(\"96\u00f78=12, 0\", \"22\u00f72=11, 0\")\n    ,@(\"83\u00f76=13, 5\", \"79\u00f73=26, 1\")\n    ,@(\"35\u00f77=5, 0\", \"83\u00f77=11, 6\")\n    ,@(\"66\u00f75=13, 1\", \"65\u00f72=32, 1\")\n    ,@(\"44\u00f79=4, 8\", \"91\u00f75=18, 1\")\n    ,@(\"94\u00f76=15, 4\", \"96\u00f77=13, 5\")\n    ,@(\"82\u00f76=13, 4\", \"85\u00f73=28, 1\")\n    ,@(\"87\u00f79=9, 6\", \"67\u00f79=7, 4\")\n    ,@(\"11\u00f74=2, 3\", \"25\u00f72=12, 1\")\n    ,@(\"38\u00f79=4, 2\", \"56\u00f72=28, 0\")\n    ,@(\"33\u00f74=8, 1\", \"25\u00f74=6, 1\")\n    ,@(\"23\u00f79=2, 5\", \"71\u00f74=17, 3\")\n    ,@(\"12\u00f79=1, 3\", \"58\u00f77=8, 2\")\n    ,@(\"44\u00f77=6, 2\", \"54\u00f73=18, 0\")\n    ,@(\"81\u00f77=11, 4\", \"41\u00f73=13, 2\")\n    ,@(\"74\u00f77=10, 4\", \"80\u00f75=16, 0\")\n    ,@(\"55\u00f78=6, 7\", \"62\u00f76=10, 2\")\n    ,@(\"32\u00f79=3, 5\", \"97\u00f75=19, 2\")\n    ,@(\"37\u00f77=5, 2\", \"69\u00f79=7, 6\")\n    ,@(\"32\u00f75=6, 2\", \"72\u00f74=18, 0\")\n    ,@(\"50\u00f73=16, 2\", \"76\u00f76=12, 4\")\n    ,@(\"38\u00f75=7, 3\", \"89\u00f79=9, 8\")\n    ,@(\"84\u00f73=28, 0\", \"14\u00f78=1, 6\")\n)\n\nfunction Replace-AllText($doc, $oldText, $newText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# Phase 1: swap each old value for a unique placeholder so a later\n# replacement never matches text an earlier one just inserted (some\n# new values are equal to other rows' old values).\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $placeholder = \"@@PLACEHOLDER\" + $i + \"@@\"\n    Replace-AllText $d $replacements[$i][0] $placeholder\n}\n\n# Phase 2: swap each placeholder for its final value.\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $placeholder = \"@@PLACEHOLDER\" + $i + \"@@\"\n    Replace-AllText $d $placeholder $replacements[$i][1]\n}\n"}
